# TN7_Test.xlsx - "working on main functions and building"
#
# Restructure the "Buildings" sheet from a per-building/per-gender single
# row with 4 floors (Rooms/Capacity pairs) into a per-building/per-gender
# header row followed by one row per floor with up to 3 room-type
# capacity/quantity pairs. Also drops the now-unused "UBC"/"Sutton"
# buildings, makes "Buildings" the active sheet/tab, and updates the
# selection on that sheet.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Buildings")

# --- Buildings sheet: rebuild the data area -------------------------------

# Drop the now-unused trailing columns (J:L) entirely - the new layout only
# uses columns A:I.
$ws2.Range("J:L").Delete()

# Clear out all the old data/header rows below row 1 so the sheet can be
# rebuilt from scratch (row 1 itself keeps its existing formatting - only
# its text changes below).
$ws2.Range("A2:I10").ClearContents()

# Update the column headers (row 1) for the new field names.
$ws2.Range("A1").Value = "Building"
$ws2.Range("B1").Value = "Gender"
$ws2.Range("C1").Value = "Number of Floors"
$ws2.Range("D1").Value = "Room Type 1 Capacity"
$ws2.Range("E1").Value = "Room Type 1 Quantity"
$ws2.Range("F1").Value = "Room Type 2 Capacity"
$ws2.Range("G1").Value = "Room Type 2 Quantity"
$ws2.Range("H1").Value = "Room Type 3 Capacity"
$ws2.Range("I1").Value = "Room Type 3 Quantity"

# Koinonia Christian / Male - 3 floors
$ws2.Range("A2").Value = "Koinonia Christian"
$ws2.Range("B2").Value = "Male"
$ws2.Range("C2").Value = 3

$ws2.Range("D3").Value = 2
$ws2.Range("E3").Value = 2
$ws2.Range("F3").Value = 3
$ws2.Range("G3").Value = 3

$ws2.Range("D4").Value = 2
$ws2.Range("E4").Value = 3
$ws2.Range("F4").Value = 3
$ws2.Range("G4").Value = 1

$ws2.Range("D5").Value = 2
$ws2.Range("E5").Value = 7

# Stratford Christian / Female - 4 floors
$ws2.Range("A6").Value = "Stratford Christian"
$ws2.Range("B6").Value = "Female"
$ws2.Range("C6").Value = 4

$ws2.Range("D7").Value = 2
$ws2.Range("E7").Value = 5
$ws2.Range("F7").Value = 3
$ws2.Range("G7").Value = 1
$ws2.Range("H7").Value = 4
$ws2.Range("I7").Value = 1

$ws2.Range("D8").Value = 2
$ws2.Range("E8").Value = 4

$ws2.Range("D9").Value = 3
$ws2.Range("E9").Value = 4

$ws2.Range("D10").Value = 2
$ws2.Range("E10").Value = 3
$ws2.Range("F10").Value = 3
$ws2.Range("G10").Value = 4
$ws2.Range("H10").Value = 4
$ws2.Range("I10").Value = 5

# --- View state: make "Buildings" the active/selected sheet+cell ---------

$ws2.Activate() | Out-Null
$ws2.Range("K5").Select() | Out-Null
